$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 474; existing rows 474:492 shift down to 476:494,
# carrying their original formatting (including the date style on column D).
$ws.Rows("474:475").Insert()

# New weekly price entries (Feria Lagunitas de Puerto Montt - Pimiento - Zafiro rojo)
$ws.Range("A474").Value = 4
$ws.Range("B474").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C474").Value = "Los Lagos"
$ws.Range("D474").Value = 44509
$ws.Range("E474").Value = 10
$ws.Range("F474").Value = 100112002
$ws.Range("G474").Value = "Pimiento"
$ws.Range("H474").Value = "Zafiro rojo"
$ws.Range("I474").Value = "Primera"
$ws.Range("J474").Value = 80
$ws.Range("K474").Value = 50000
$ws.Range("L474").Value = 50000
$ws.Range("M474").Value = 50000
$ws.Range("N474").Value = "$/caja 15 kilos"
$ws.Range("O474").Value = "Región de Arica y Parinacota"
$ws.Range("P474").Value = 3333
$ws.Range("Q474").Value = 15
$ws.Range("R474").Value = "Hortaliza"

$ws.Range("A475").Value = 4
$ws.Range("B475").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C475").Value = "Los Lagos"
$ws.Range("D475").Value = 44509
$ws.Range("E475").Value = 10
$ws.Range("F475").Value = 100112002
$ws.Range("G475").Value = "Pimiento"
$ws.Range("H475").Value = "Zafiro rojo"
$ws.Range("I475").Value = "Segunda"
$ws.Range("J475").Value = 80
$ws.Range("K475").Value = 47000
$ws.Range("L475").Value = 47000
$ws.Range("M475").Value = 47000
$ws.Range("N475").Value = "$/caja 15 kilos"
$ws.Range("O475").Value = "Región de Arica y Parinacota"
$ws.Range("P475").Value = 3133
$ws.Range("Q475").Value = 15
$ws.Range("R475").Value = "Hortaliza"
